# Generate Report for Handback
# The handback files are now in sync with en-US, so the localization-status
# report is refreshed: status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", the handback timestamps advance, and the
# stale "version not latest" error details are cleared now that the issue is
# resolved.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.144371396019366
$wsOverview.Columns.Item(6).ColumnWidth = 29.144371396019366

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("K2").Value = "2016-09-03 16:56:28"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsZhCn.Columns.Item(16).ColumnWidth = 12.913719540550566

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("K2").Value = "2016-09-03 16:56:35"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsDeDe.Columns.Item(16).ColumnWidth = 12.913719540550566
